$wb = $excel.ActiveWorkbook

# --- 1. Add the new "userInfo" worksheet after the existing "user" sheet ---
$userSheet = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "userInfo"

# --- 2. Header row (row 1) ---
$ws.Range("A1").Value = "testcase_name"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "lastname"
$ws.Range("F1").Value = "password"
$ws.Range("G1").Value = "confirmpassword"
$ws.Range("H1").Value = "role"
$ws.Range("I1").Value = "isAdmin"
$ws.Range("J1").Value = "admintext"
$ws.Range("K1").Value = "ReportName"
$ws.Range("L1").Value = "ReportFolder"
$ws.Range("M1").Value = "PrimaryModule"
$ws.Range("N1").Value = "RelatedModule"
$ws.Range("O1").Value = "ConditionField"
$ws.Range("P1").Value = "CompareField"
$ws.Range("Q1").Value = "GroupField"

# Highlight the header row with a yellow fill
$ws.Range("A1:Q1").Interior.Color = 65535

# --- 3. Data row (row 2) ---
$ws.Range("A2").Value = "create user with admin checked"
$ws.Range("B2").Value = "TestDerek O'Connell"
$ws.Range("C2").Value = "TestDerek OConnell"
$ws.Range("D2").Value = "adminUser@gmail.com"
$ws.Range("E2").Value = "admin"
$ws.Range("F2").Value = "Admin@12345"
$ws.Range("G2").Value = "Admin@12345"
$ws.Range("H2").Value = "IT"
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = "Yes"
$ws.Range("K2").Value = "Chart_Report"
$ws.Range("L2").Value = "Test Chethana"
$ws.Range("M2").Value = "Contacts"
$ws.Range("N2").Value = "Calendar"
$ws.Range("O2").Value = "First Name"
$ws.Range("P2").Value = "Test"
$ws.Range("Q2").Value = "First Name"

# --- 4. Selections as in the target workbook ---
# "user" sheet: selection becomes a range, no longer the active tab
$userSheet.Range("A1:I2").Select() | Out-Null

# "userInfo" sheet: becomes the active tab with C2 selected
$ws.Range("C2").Select() | Out-Null
